{"js": "// Regenerate the lattice-multiplication practice table: every cell keeps\n// its original 5-line layout (problem / digits / divider / two partial\n// rows) but gets new numbers, matching a freshly generated worksheet.\n// Word represents the `<w:br/>` line breaks inside a cell as the\n// vertical-tab character (\\u000b) when read/written through Range.text,\n// so each cell's new content is built as those 5 lines joined by \\u000b.\n\nconst NL = \"\\u000b\";\n\n// New cell text, in document order (row-major, same 5x3 shape as before).\nconst newCells = [\n  [\"82 x 94\", \"  9    4\", \"  ----\", \"8|    |\", \"2|    |\"],\n  [\"80 x 24\", \"  2    4\", \"  ----\", \"8|    |\", \"0|    |\"],\n  [\"82 x 27\", \"  2    7\", \"  ----\", \"8|    |\", \"2|    |\"],\n  [\"23 x 78\", \"  7    8\", \"  ----\", \"2|    |\", \"3|    |\"],\n  [\"43 x 32\", \"  3    2\", \"  ----\", \"4|    |\", \"3|    |\"],\n  [\"81 x 60\", \"  6    0\", \"  ----\", \"8|    |\", \"1|    |\"],\n  [\"83 x 35\", \"  3    5\", \"  ----\", \"8|    |\", \"3|    |\"],\n  [\"17 x 38\", \"  3    8\", \"  ----\", \"1|    |\", \"7|    |\"],\n  [\"41 x 27\", \"  2    7\", \"  ----\", \"4|    |\", \"1|    |\"],\n  [\"98 x 57\", \"  5    7\", \"  ----\", \"9|    |\", \"8|    |\"],\n  [\"53 x 38\", \"  3    8\", \"  ----\", \"5|    |\", \"3|    |\"],\n  [\"42 x 11\", \"  1    1\", \"  ----\", \"4|    |\", \"2|    |\"],\n  [\"74 x 58\", \"  5    8\", \"  ----\", \"7|    |\", \"4|    |\"],\n  [\"54 x 50\", \"  5    0\", \"  ----\", \"5|    |\", \"4|    |\"],\n  [\"24 x 67\", \"  6    7\", \"  ----\", \"2|    |\", \"4|    |\"],\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = 3;\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    if (idx >= newCells.length) break;\n    const cell = table.getCell(r, c);\n    const text = newCells[idx].join(NL);\n    cell.getRange().insertText(text, Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Regenerate the lattice-multiplication practice table: every cell keeps\n# its original 5-line layout (problem / digits / divider / two partial\n# rows) but gets new numbers, matching a freshly generated worksheet.\n# Word exposes the `<w:br/>` line breaks inside a cell as the\n# vertical-tab character (Chr(11)) through Range.Text, so each cell's new\n# content is built as those 5 lines joined by Chr(11).\n\n$nl = [char]11\n\n# New cell text, in document order (row-major, same 5x3 shape as before).\n$newCells = @(\n    @(\"82 x 94\", \"  9    4\", \"  ----\", \"8|    |\", \"2|    |\"),\n    @(\"80 x 24\", \"  2    4\", \"  ----\", \"8|    |\", \"0|    |\"),\n    @(\"82 x 27\", \"  2    7\", \"  ----\", \"8|    |\", \"2|    |\"),\n    @(\"23 x 78\", \"  7    8\", \"  ----\", \"2|    |\", \"3|    |\"),\n    @(\"43 x 32\", \"  3    2\", \"  ----\", \"4|    |\", \"3|    |\"),\n    @(\"81 x 60\", \"  6    0\", \"  ----\", \"8|    |\", \"1|    |\"),\n    @(\"83 x 35\", \"  3    5\", \"  ----\", \"8|    |\", \"3|    |\"),\n    @(\"17 x 38\", \"  3    8\", \"  ----\", \"1|    |\", \"7|    |\"),\n    @(\"41 x 27\", \"  2    7\", \"  ----\", \"4|    |\", \"1|    |\"),\n    @(\"98 x 57\", \"  5    7\", \"  ----\", \"9|    |\", \"8|    |\"),\n    @(\"53 x 38\", \"  3    8\", \"  ----\", \"5|    |\", \"3|    |\"),\n    @(\"42 x 11\", \"  1    1\", \"  ----\", \"4|    |\", \"2|    |\"),\n    @(\"74 x 58\", \"  5    8\", \"  ----\", \"7|    |\", \"4|    |\"),\n    @(\"54 x 50\", \"  5    0\", \"  ----\", \"5|    |\", \"4|    |\"),\n    @(\"24 x 67\", \"  6    7\", \"  ----\", \"2|    |\", \"4|    |\")\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -ge $newCells.Count) { continue }\n        $lines = $newCells[$idx]\n        $text = [string]::Join($nl, $lines)\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $text\n        $idx++\n    }\n}\n"}
